$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.395.55"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "3.145.11"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'604.16"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "'150.09"
$ws.Range("E6").Value = "  -4.53%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "3.142.92"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("E10").Value = "  -4.41%  "
$ws.Range("D11").Value = "'5.63"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("E12").Value = "  -4.79%  "
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("D14").Value = "'37.03"
$ws.Range("E14").Value = "  -4.28%  "
$ws.Range("D15").Value = "3.620.64"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "64.446.12"
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "3.146.61"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("D20").Value = "'484.45"
$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").Value = "'13.91"
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").Value = "'84.51"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'2.93"
$ws.Range("E27").Value = "  -2.42%  "
$ws.Range("D28").Value = "'8.69"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'7.12"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("D32").Value = "'2.73"
$ws.Range("E32").Value = "  -7.03%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'26.90"
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("D36").Value = "'6.13"
$ws.Range("E36").Value = "  -5.97%  "
$ws.Range("D37").Value = "'3.29"
$ws.Range("E37").Value = "  +7.11%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "0.0₃0760"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'451.43"
$ws.Range("E40").Value = "  -10.56%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.125"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0403"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.898.39"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.42"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -7.66%  "
$ws.Range("D47").Value = "'27.01"
$ws.Range("E47").Value = "  -4.35%  "
$ws.Range("D48").Value = "'0.998"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").Value = "'120.25"
$ws.Range("E51").Value = "  -1.45%  "

# Reset style on cells that were force-typed as text via leading quote,
# so no stray quotePrefix style index is left attached to the cell.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
